# Add the missing "Events" row data that was causing a KeyError on page 1 when predicting.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

$ws.Cells.Item(2, 1).Value = 77
$ws.Cells.Item(2, 2).Value = "Event reporting 4"
$ws.Cells.Item(2, 3).Value = 222.43077994188189
$ws.Cells.Item(2, 4).Value = "Thrombotic event not related to the area intervened"

$ws.Range("H7").Select()
